$d = $word.ActiveDocument

function Find-ParagraphIndex($predicateText) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $t = $paras.Item($i).Range.Text
        if ($t.StartsWith($predicateText)) {
            return $i
        }
    }
    throw "Paragraph starting with '$predicateText' not found"
}

function Replace-ParagraphXml($paraIndex, $xmlFrag) {
    # Re-fetch the paragraph Range fresh (not reusing any Range that has
    # already been touched by Find or a previous InsertXML call) - doing so
    # keeps InsertXML's content-replacement behaviour well defined, in
    # particular when the paragraph contains hyperlink fields.
    $para = $d.Paragraphs.Item($paraIndex)
    $r = $para.Range
    $r.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# Change 1: paragraph "Schritte um Projekt auf Maschine zum Laufen zu
# bringen:" - originally split across two runs with a proofErr gramStart/
# gramEnd pair around "Schritte". Merge into a single run and drop the
# proofErr markers.
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphIndex "Schritte um Projekt"
$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="071FE054" w14:textId="77777777" w:rsidR="00626276" w:rsidRPr="00626276" w:rsidRDefault="00626276" w:rsidP="00626276">
            <w:pPr>
              <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              <w:spacing w:after="240" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
              <w:t>Schritte um Projekt auf Maschine zum Laufen zu bringen:</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
Replace-ParagraphXml $idx1 $xml1

# ---------------------------------------------------------------------------
# Change 2: simple text replacement, no structural changes needed.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Dateien unter 'htdocs' im XAMPP-Ordner speichern",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dateien im XAMPP-Ordner unter C:\xampp\htdocs in einem neuen Ordner namens HelpDesk speichern", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: paragraph "Im Browser 'http://localhost/HelpDesk/' aufrufen um
# Ticket abschicken zu können" - merge the trailing runs split around the
# proofErr gramStart/gramEnd pair on "aufrufen" into a single run, keeping
# the leading run and the hyperlink run untouched.
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphIndex "Im Browser 'http://localhost/HelpDesk/' aufrufen"
$xml3 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p w14:paraId="47504C2C" w14:textId="77777777" w:rsidR="00626276" w:rsidRPr="00626276" w:rsidRDefault="00626276" w:rsidP="00626276">
            <w:pPr>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              <w:spacing w:before="60" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="00626276">
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
              <w:t>Im Browser '</w:t>
            </w:r>
            <w:hyperlink r:id="rId6" w:history="1">
              <w:r w:rsidRPr="00626276">
                <w:rPr>
                  <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                  <w:color w:val="0000FF"/>
                  <w:sz w:val="24"/>
                  <w:szCs w:val="24"/>
                  <w:u w:val="single"/>
                  <w:lang w:eastAsia="de-DE"/>
                </w:rPr>
                <w:t>http://localhost/HelpDesk/</w:t>
              </w:r>
            </w:hyperlink>
            <w:r w:rsidRPr="00626276">
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
              <w:t>' aufrufen um Ticket abschicken zu können</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
Replace-ParagraphXml $idx3 $xml3

# ---------------------------------------------------------------------------
# Change 4: paragraph "Im Browser 'http://localhost/HelpDesk/ticketsystem.php'
# aufrufen um Ticketübersicht zu sehen und Ticket auf Status gelöst setzen zu
# können." - same kind of run merge / proofErr removal as change 3.
# ---------------------------------------------------------------------------
$idx4 = Find-ParagraphIndex "Im Browser 'http://localhost/HelpDesk/ticketsystem.php'"
$xml4 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p w14:paraId="6E87073A" w14:textId="77777777" w:rsidR="00626276" w:rsidRPr="00626276" w:rsidRDefault="00626276" w:rsidP="00626276">
            <w:pPr>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              <w:spacing w:before="60" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="00626276">
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
              <w:t>Im Browser '</w:t>
            </w:r>
            <w:hyperlink r:id="rId7" w:history="1">
              <w:r w:rsidRPr="00626276">
                <w:rPr>
                  <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                  <w:color w:val="0000FF"/>
                  <w:sz w:val="24"/>
                  <w:szCs w:val="24"/>
                  <w:u w:val="single"/>
                  <w:lang w:eastAsia="de-DE"/>
                </w:rPr>
                <w:t>http://localhost/HelpDesk/ticketsystem.php</w:t>
              </w:r>
            </w:hyperlink>
            <w:r w:rsidRPr="00626276">
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="24292F"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
              <w:t>' aufrufen um Ticketübersicht zu sehen und Ticket auf Status gelöst setzen zu können.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
Replace-ParagraphXml $idx4 $xml4
